# Apply the "Committing the new tests and changes in the tests." edit to
# Sheet1 of the workbook: rename/clear several Genetic-Feature-with-evidence
# protein/sequence columns, shorten the big DNA sequence blob, and move the
# current selection / view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) -------------------------------------------------------
# T1 "proteinSymbol" -> "Symbol"; U1..Z1 and AB1 cleared out entirely.
$ws.Range("T1").Value = "Symbol"
$ws.Range("U1").ClearContents()
$ws.Range("V1").ClearContents()
$ws.Range("W1").ClearContents()
$ws.Range("X1").ClearContents()
$ws.Range("Y1").ClearContents()
$ws.Range("Z1").ClearContents()
$ws.Range("AB1").ClearContents()

# --- Row 2 (data) -----------------------------------------------------------
# L2 held a huge DNA sequence string; trim it down to just its tail.
$ws.Range("L2").Value = "
amakpgqndk lrhagiidiq  fqrvpcnhpg lnvnfqverg 181 snpnylavlv efanregtvv qmdlmesrng rptgywtamr hswgaiwrmd srrrlqgpfs 241 lrirsesgkt lvakqvipan wrpdtnyrsn vqfr"

$ws.Range("R2").Value = "VVNNAAAA"
$ws.Range("T2").Value = "selenium_GF1"
$ws.Range("U2").ClearContents()
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Column widths ----------------------------------------------------------
$ws.Range("T1").EntireColumn.ColumnWidth = 28
$ws.Range("W1").EntireColumn.ColumnWidth = 18
$ws.Range("X1").EntireColumn.ColumnWidth = 23.140625
$ws.Range("Z1").EntireColumn.ColumnWidth = 27.28515625
$ws.Cells.Item(1, 53).EntireColumn.ColumnWidth = 24

# --- Page setup / view -------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 10

$ws.Range("Q3").Select()
